# Natmi following Dr Hou advice
# Update the existing Fgf15-Fgfr3 LR-pair rows (2 and 3) with recomputed
# statistics and introduce the new "FAPs" cluster, adding rows 4-7 for the
# additional sending/target cluster combinations that now appear in the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf15"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.029478
$ws.Range("H2").Value = 0.088434
$ws.Range("I2").Value = 0.1535387136874709
$ws.Range("J2").Value = 0.1535387136874709
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.243417666666667
$ws.Range("N2").Value = 15.730253
$ws.Range("O2").Value = 0.8253998362974575
$ws.Range("P2").Value = 0.8253998362974574
$ws.Range("Q2").Value = 0.154565465978
$ws.Range("R2").Value = 1.391089193802
$ws.Range("S2").Value = 0.1267308291429607
$ws.Range("T2").Value = 0.1267308291429607

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf15"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.029478
$ws.Range("H3").Value = 0.088434
$ws.Range("I3").Value = 0.1535387136874709
$ws.Range("J3").Value = 0.1535387136874709
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.6792986666666666
$ws.Range("N3").Value = 2.037896
$ws.Range("O3").Value = 0.1069327381315001
$ws.Range("P3").Value = 0.1069327381315001
$ws.Range("Q3").Value = 0.020024366096
$ws.Range("R3").Value = 0.180219294864
$ws.Range("S3").Value = 0.01641831506378969
$ws.Range("T3").Value = 0.01641831506378969

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf15"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.029478
$ws.Range("H4").Value = 0.088434
$ws.Range("I4").Value = 0.1535387136874709
$ws.Range("J4").Value = 0.1535387136874709
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4298626666666667
$ws.Range("N4").Value = 1.289588
$ws.Range("O4").Value = 0.06766742557104236
$ws.Range("P4").Value = 0.06766742557104236
$ws.Range("Q4").Value = 0.012671491688
$ws.Range("R4").Value = 0.114043425192
$ws.Range("S4").Value = 0.01038956948072052
$ws.Range("T4").Value = 0.01038956948072052

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Fgf15"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1625126666666667
$ws.Range("H5").Value = 0.487538
$ws.Range("I5").Value = 0.846461286312529
$ws.Range("J5").Value = 0.846461286312529
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.243417666666667
$ws.Range("N5").Value = 15.730253
$ws.Range("O5").Value = 0.8253998362974575
$ws.Range("P5").Value = 0.8253998362974574
$ws.Range("Q5").Value = 0.8521217874571111
$ws.Range("R5").Value = 7.669096087114
$ws.Range("S5").Value = 0.6986690071544968
$ws.Range("T5").Value = 0.6986690071544966

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Fgf15"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1625126666666667
$ws.Range("H6").Value = 0.487538
$ws.Range("I6").Value = 0.846461286312529
$ws.Range("J6").Value = 0.846461286312529
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.6792986666666666
$ws.Range("N6").Value = 2.037896
$ws.Range("O6").Value = 0.1069327381315001
$ws.Range("P6").Value = 0.1069327381315001
$ws.Range("Q6").Value = 0.1103946377831111
$ws.Range("R6").Value = 0.9935517400480001
$ws.Range("S6").Value = 0.09051442306771036
$ws.Range("T6").Value = 0.09051442306771036

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Fgf15"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1625126666666667
$ws.Range("H7").Value = 0.487538
$ws.Range("I7").Value = 0.846461286312529
$ws.Range("J7").Value = 0.846461286312529
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4298626666666667
$ws.Range("N7").Value = 1.289588
$ws.Range("O7").Value = 0.06766742557104236
$ws.Range("P7").Value = 0.06766742557104236
$ws.Range("Q7").Value = 0.06985812826044445
$ws.Range("R7").Value = 0.6287231543440002
$ws.Range("S7").Value = 0.05727785609032184
$ws.Range("T7").Value = 0.05727785609032184

